$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---
# Overview sheet keeps one status column per locale (E = zh-cn, F = de-de)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Per-locale detail sheets keep the status in column C
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the now-shorter Status column(s) to fit the new text ---
# New text is shorter than the old, so the report-generation step re-sized
# the Status columns down. Closest attainable width via the Excel
# character-width (pixel-quantized) model for the recorded target width.
$newStatusWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth
